$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new "season record" columns (AD:AF), reusing the
# existing header style (s="1") from the adjacent header cell AC1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row (2-37).
$lastRow = 37
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 86
    $ws.Cells.Item($r, 31).Value = 76
    $ws.Cells.Item($r, 32).Value = 0
}
